# Project idea presentation - remove the trailing period after the final
# ”Vote” entity. reference in the last bullet of the body placeholder on
# slide 1 ("... in the "Vote" entity." -> "... in the "Vote" entity").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the body placeholder shape that holds the project description text
# (rather than assuming a fixed shape index).
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text -like "*entity.*") {
            $target = $candidate
        }
    }
}

$tr = $target.TextFrame.TextRange
$fullText = $tr.Text

# Find the very last occurrence of the word "entity." (with its trailing
# full stop) - this is the sentence-ending reference to the weak "Vote"
# entity at the end of the 4th paragraph.
$needle = "entity."
$idx0 = $fullText.LastIndexOf($needle)

if ($idx0 -ge 0) {
    # Convert the 0-based .NET string index of the period to PowerPoint's
    # 1-based TextRange character position.
    $periodPos = $idx0 + $needle.Length
    $periodRange = $tr.Characters($periodPos, 1)
    if ($periodRange.Text -eq ".") {
        $periodRange.Delete()
    }
}

Write-Output "Updated text: $($target.TextFrame.TextRange.Text)"
